$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-12 from 45174 (2023-09-05)
# to 45175 (2023-09-06).
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
